$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "AINA+ 1"
$ws.Range("B3").Value  = "AINA– 2"
$ws.Range("B4").Value  = "REFHA 3"
$ws.Range("B5").Value  = "REFHA 4"
$ws.Range("B6").Value  = "REFLA 5"
$ws.Range("B7").Value  = "VDD 7"
$ws.Range("B8").Value  = "VDD 7"
$ws.Range("B9").Value  = "CLKA 8"
$ws.Range("B10").Value = "CLKB 9"
$ws.Range("B11").Value = "VDD 10"
$ws.Range("B12").Value = "REFLB 11"
$ws.Range("B13").Value = "REFLB 12"
$ws.Range("B14").Value = "REFHB 13"
$ws.Range("B15").Value = "REFHB 14"
$ws.Range("B16").Value = "AINB– 15"
$ws.Range("B17").Value = "AINB+ 16"
$ws.Range("B18").Value = "GND 17"
$ws.Range("B19").Value = "VDD 18"
$ws.Range("B20").Value = "VDD 18"
$ws.Range("B21").Value = "VCMB 20"
$ws.Range("B22").Value = "MUX 21"
$ws.Range("B23").Value = "OEB 23"
$ws.Range("B24").Value = "OEB 23"
$ws.Range("B25").Value = "NC 24"
$ws.Range("B26").Value = "NC 25"
$ws.Range("B27").Value = "DB0 26"
$ws.Range("B28").Value = "DB1 27"
$ws.Range("B29").Value = "DB2 28"
$ws.Range("B30").Value = "DB3 29"
$ws.Range("B31").Value = "DB4 30"
$ws.Range("B32").Value = "OGND 31"
$ws.Range("B33").Value = "OVDD 32"
$ws.Range("B34").Value = "33 DB5"
$ws.Range("B35").Value = "34 DB6"
$ws.Range("B36").Value = "35 DB7"
$ws.Range("B37").Value = "36 DB8"
$ws.Range("B38").Value = "37 DB9"
$ws.Range("B39").Value = "37 DB9"
$ws.Range("B40").Value = "39 DB11"
$ws.Range("B41").Value = "40 OFB "
$ws.Range("B42").Value = "41 NC"
$ws.Range("B43").Value = "42 NC"
$ws.Range("B44").Value = "43 DA0"
$ws.Range("B45").Value = "44 DA1 "
$ws.Range("B46").Value = "45 DA2"
$ws.Range("B47").Value = "46 DA3"
$ws.Range("B48").Value = "47 DA4 "
$ws.Range("B49").Value = "48 DA5 "
$ws.Range("B50").Value = "49 OVDD"
$ws.Range("B51").Value = "50 OGND"
$ws.Range("B52").Value = "51 DA6"
$ws.Range("B53").Value = "52 DA7"
$ws.Range("B54").Value = "53 DA8"
$ws.Range("B55").Value = "54 DA9"
$ws.Range("B56").Value = "55 DA10"
$ws.Range("B57").Value = "56 DA11"
$ws.Range("B58").Value = "57 OFA"
$ws.Range("B59").Value = "58 OEA"
$ws.Range("B60").Value = "59 SHDNA"
$ws.Range("B61").Value = "60 MODE"
$ws.Range("B62").Value = "61 VCMA"
$ws.Range("B63").Value = "62 SENSEA"
$ws.Range("B64").Value = "63 VDD"
$ws.Range("B65").Value = "64 GND"
$ws.Range("B66").Value = "57 OFA"
